# Add older catch limits to the Landings_FY1998 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 13 (e.g. species row 10): add Percent of TAC and Target TAC values
$ws.Range("B13").Value = 12110
$ws.Range("D13").Value = 55045.454545454544

# Row 15 (species row 11): add Percent of TAC and Target TAC values
$ws.Range("B15").Value = 758
$ws.Range("D15").Value = 75800

# Row 17 (species row 12): add Target TAC value
$ws.Range("D17").Value = 56550

# Row 19 (species row 13): add Target TAC value
# Also align its formatting with the other "Target TAC" cells (copy format from D17)
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("D19").Value = 55621.428571428565

# Row 25 (species row 15): add Target TAC value
$ws.Range("D25").Value = 57828.57142857142

# Update the active selection to reflect where editing left off
$ws.Range("L17").Select() | Out-Null
